$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

$data = @(
    @('E0032', 'Indicating Controller identity and details of representative', 'This example describes a controller with details for name, address, and contact, and its representatives for EU and US jurisdictions', 'E0032.ttl', 'dpv:DataController,dpv:Representative'),
    @('E0033', 'Indicating Processor as the implementing entity in a process', 'This example shows a process with its controller and processor, and the explicit statement for the processor doing the collection and storing of personal data', 'E0033.ttl', 'dpv:DataController,dpv:DataProcessor'),
    @('E0034', 'Specifying recipients of data', 'This example first shows a process where the processor is denoted as the data recipient, and then another process where a third party is denoted as the recipient. A joint controller agreement is then shown which specifies that only one of the controllers involved is a data recipient.', 'E0034.ttl', 'dpv:DataRecipient'),
    @('E0035', 'Specifying data exporters and importers', 'This example shows how data exporters and importers for a transfer can be indicated using DPV. It also shows how the locations associated with exporters and importers are relevant to determine cross-border transfers.', 'E0035.ttl', 'dpv:DataImporter,dpv:DataExporter,dpv:Transfer'),
    @('E0036', 'Indicate relevant authority for processing', 'This example shows how a DPA can be associated with processing, and the use of LEGAL extensions to obtain DPA information. It also shows how DPAs can be ''discovered'' by using the location (jurisdiction) and applicable law concepts.', 'E0036.ttl', 'dpv:DataProtectionAuthority'),
    @('E0037', 'Indicating type of organisation and involvement of specific orgnisational units', 'This example involves an organisation that is a NGO, and that it has Marketing, HR, and IT departments. The HR and IT departments are responsible for specific processes.', 'E0037.ttl', 'dpv:OrganisationalUnit,dpv:isImplementedByEntity'),
    @('E0038', 'Indicating subsidiaries of an organisation', 'This example shows the existence of two subsidiaries associated with an organisation and their locations ', 'E0038.ttl', 'dpv:Subsidiary,dpv:hasSubsidiary,dpv:isSubsidiaryFor'),
    @('E0039', 'Indicating involvement of data subjects', 'This example shows the different ways in which data subjects can be indicated in a process. It also shows how the DPV taxonomy of data subjects can be used to combine concepts to accurately represent the data subject involved. And it also shows how information associated with specific data subjects such as identifiers can be indicated.', 'E0039.ttl', 'dpv:DataSubject'),
)

$startRow = 33
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = "ttl"
    $ws.Cells.Item($r, 6).Value = "file"
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 9).Value = "accepted"
    $ws.Cells.Item($r, 10).Value = (Get-Date -Year 2024 -Month 6 -Day 10)
    $ws.Cells.Item($r, 11).Value = "Harshvardhan J. Pandit"
}
